$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.507985234260559
$ws.Range("B1").Value = 1.530914187431335
$ws.Range("C1").Value = 1.571980953216553
$ws.Range("D1").Value = 2.081084489822388
$ws.Range("E1").Value = 3.52023983001709
